$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill row 7 (powerUps glitch entry)
$ws.Range("A7").Value = "powerUps"
$ws.Range("B7").Value = "hitting brick with ball and hitting powerup with paddle"
$ws.Range("C7").Value = "When ball breaks brick a chance of powerup spawning should drop down and player gets the powerrup if they catch it with paddle"
$ws.Range("D7").Value = "the power up spawns and drops down but when the player hits it sometimes runs twice changing powerupcount down by 2 causing glithces"
$ws.Range("E7").Value = "the fucniton was getting called once but was happening twice for some reason because of for each loop fixed glithc by moving for loop for ball after the pUpCount is lowered"

# Fill row 8 (Paddle moving glitch entry)
$ws.Range("A8").Value = "Paddle moveing"
$ws.Range("B8").Value = "Mouse movement"
$ws.Range("D8").Value = "THe paddle glitches out when mouse is of the canvas"
$ws.Range("E8").Value = "this is because the mousex was undefined. I fixed this by making it only change paddle if mouse x is with in the canvas and go to the middle at the start of the game"
$ws.Range("C8").Value = "The paddle should follow the mouse horizontaly"

$ws.Rows.Item(7).RowHeight = 66
$ws.Rows.Item(8).RowHeight = 66

# Update the view/selection
$ws.Range("D7").Select()
